$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Configuration")

# Remove the "MARKETS" row (row 5) entirely - shifts all rows below it up by one.
$ws.Rows.Item(5).Delete() | Out-Null

# Update the remaining values to their new figures (rows renumbered after the delete).
$ws.Range("B1").Value = 10    # PERIODS
$ws.Range("B2").Value = 2     # AGENTS
$ws.Range("B3").Value = 0     # CONTACTS
# B4 (FRIENDS) and B5 (LEVELS) are unchanged.
$ws.Range("B6").Value = 0     # REPETITIONS
$ws.Range("B7").Value = 1     # GUI
# B8 (BASE), B9 (MEMORY) and B10 (SAVED_ENDORSEMENTS) are unchanged.

# Move the active selection to B4, matching the saved view state.
$ws.Range("B4").Select() | Out-Null
